$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("Sequences")
$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "CCTGGCCCTCTCCAACGTCAAAGGGCGACTTGACGGGGAAAGttCCTACTTCATCCATTAAATCC"
$arr[0,1] = "CATAAAGTGGTGAGACGGGCAACAGCTGAGAAAGCGAAAGGAttCCTACTTCATCCATTAAATCC"
$arr[0,2] = "TCTGTGGTGGCTCACAATTCCACACAACCGGTCACGCTGCGCttCCTACTTCATCCATTAAATCC"
$arr[0,3] = "ATCCGCCGGGATCCAGCGCAGTGTCACTCGCCGCTACAGGGCttCCTACTTCATCCATTAAATCC"
$arr[0,4] = "TTTTCGTCTTCAGCGGGGTCATTGCAGGTATAACGTGCTTTCttCCTACTTCATCCATTAAATCC"
$arr[0,5] = "TTTCTCCGTTGCTGATTGCCGTTCCGGCAGGAGGCCGATTAAttCCTACTTCATCCATTAAATCC"
$arr[0,6] = "GCTATTACGGTTTACCAGTCCCGGAATTGAATCCTGAGAAGTttCCTACTTCATCCATTAAATCC"
$arr[0,7] = "GACCGTAATCTGTTGGGAAGGGCGATCGAAAAGAGTCTGTCCttCCTACTTCATCCATTAAATCC"
$arr[0,8] = "GGAAGATTGCGTCGGATTCTCCGTGGGACTTCTTTGATTAGTttCCTACTTCATCCATTAAATCC"
$arr[0,9] = "AGACAGTCACCCCGGTTGATAATCAGAACTCAAACTATCGGCttCCTACTTCATCCATTAAATCC"
$arr[0,10] = "GGCAAGGCATAGGTAAAGATTCAAAAGGCCGCCAGCCATTGCttCCTACTTCATCCATTAAATCC"
$arr[0,11] = "AATATGCAATAGTAGTAGCATTAACATCACATTTTGACGCTCttCCTACTTCATCCATTAAATCC"
$arr[0,12] = "AGCGGATTGGCTGAATATAATGCTGTAGGCAGATTCACCAGTttCCTACTTCATCCATTAAATCC"
$arr[0,13] = "GATAAAAACGGTCTTTACCCTGACTATTTGGCCAACAGAGATttCCTACTTCATCCATTAAATCC"
$arr[0,14] = "TAAGAACTGCAACACTATCATAACCCTCATACGTGGCACAGAttCCTACTTCATCCATTAAATCC"
$arr[0,15] = "ACGGTGTACACTTTAATCATTGTGAATTAATGCGCGAACTGAttCCTACTTCATCCATTAAATCC"
$ws.Range("B8:Q8").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "ACGAAAGAGCCGAACTGACCAACTTTGACCGAACGAACCACCttCCTACTTCATCCATTAAATCC"
$arr[0,1] = "CGCCCACGCTACGTAATGCCACTACGAATCAGTATTAACACCttCCTACTTCATCCATTAAATCC"
$arr[0,2] = "TTTGCTAAAACCGATAGTTGCGCCGACAGCAGCAAATGAAAAttCCTACTTCATCCATTAAATCC"
$arr[0,3] = "CCACCCTCACAGACGTTAGTAAATGAATAATATCAAACCCTCttCCTACTTCATCCATTAAATCC"
$arr[0,4] = "ACAGTTAATTCAGGAGGTTTAGTACCGCACAGTTGAAAGGAAttCCTACTTCATCCATTAAATCC"
$arr[0,5] = "CCGCCGCCAGGGTCAGTGCCTTGAGTAAGGAGCACTAACAACttCCTACTTCATCCATTAAATCC"
$arr[0,6] = "GAATCAAGTCCCTCAGAGCCGCCACCAGACATTTGAGGATTTttCCTACTTCATCCATTAAATCC"
$arr[0,7] = "CCAAAGACACCATCGATAGCAGCACCGTACAACTCGTATTAAttCCTACTTCATCCATTAAATCC"
$arr[0,8] = "GAACAAAGTCAATCAATAGAAAATTCATAAAGTTTGAGTAACttCCTACTTCATCCATTAAATCC"
$arr[0,9] = "CTTTACAGAGAAGCCCTTTTTAAGAAAACCAGAAGGAGCGGAttCCTACTTCATCCATTAAATCC"
$arr[0,10] = "AACCTCCCGTTTTTGTTTAACGTCAAAAGATGGCAATTCATCttCCTACTTCATCCATTAAATCC"
$arr[0,11] = "CATCCTAATTCCGGTATTCTAAGAACGCTTCTGAATAATGGAttCCTACTTCATCCATTAAATCC"
$arr[0,12] = "GTAGGGCTTATAGATAAGTCCTGAACAATTTGCACGTAAAACttCCTACTTCATCCATTAAATCC"
$arr[0,13] = "ATCGCAAGAAAATTCTTACCAGTATAAAGGTTTAACGTCAGAttCCTACTTCATCCATTAAATCC"
$arr[0,14] = "ATATATGTGATATAACTATATGTAAATGTCGGGAGAAACAATttCCTACTTCATCCATTAAATCC"
$arr[0,15] = "TCTGTAAATTAACAATTTCATTTTTTTAATGGAAACAAGTTACAAAATCttCCTACTTCATCCATTAAATCC"
$ws.Range("B9:Q9").Value = $arr


$ws = $wb.Worksheets.Item("Descriptions")
$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "cc6hb h5 position 1 staple with tt linker and anti-Quimby handle"
$arr[0,1] = "cc6hb h5 position 2 staple with tt linker and anti-Quimby handle"
$arr[0,2] = "cc6hb h5 position 3 staple with tt linker and anti-Quimby handle"
$arr[0,3] = "cc6hb h5 position 4 staple with tt linker and anti-Quimby handle"
$arr[0,4] = "cc6hb h5 position 5 staple with tt linker and anti-Quimby handle"
$arr[0,5] = "cc6hb h5 position 6 staple with tt linker and anti-Quimby handle"
$arr[0,6] = "cc6hb h5 position 7 staple with tt linker and anti-Quimby handle"
$arr[0,7] = "cc6hb h5 position 8 staple with tt linker and anti-Quimby handle"
$arr[0,8] = "cc6hb h5 position 9 staple with tt linker and anti-Quimby handle"
$arr[0,9] = "cc6hb h5 position 10 staple with tt linker and anti-Quimby handle"
$arr[0,10] = "cc6hb h5 position 11 staple with tt linker and anti-Quimby handle"
$arr[0,11] = "cc6hb h5 position 12 staple with tt linker and anti-Quimby handle"
$arr[0,12] = "cc6hb h5 position 13 staple with tt linker and anti-Quimby handle"
$arr[0,13] = "cc6hb h5 position 14 staple with tt linker and anti-Quimby handle"
$arr[0,14] = "cc6hb h5 position 15 staple with tt linker and anti-Quimby handle"
$arr[0,15] = "cc6hb h5 position 16 staple with tt linker and anti-Quimby handle"
$ws.Range("B8:Q8").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "cc6hb h5 position 17 staple with tt linker and anti-Quimby handle"
$arr[0,1] = "cc6hb h5 position 18 staple with tt linker and anti-Quimby handle"
$arr[0,2] = "cc6hb h5 position 19 staple with tt linker and anti-Quimby handle"
$arr[0,3] = "cc6hb h5 position 20 staple with tt linker and anti-Quimby handle"
$arr[0,4] = "cc6hb h5 position 21 staple with tt linker and anti-Quimby handle"
$arr[0,5] = "cc6hb h5 position 22 staple with tt linker and anti-Quimby handle"
$arr[0,6] = "cc6hb h5 position 23 staple with tt linker and anti-Quimby handle"
$arr[0,7] = "cc6hb h5 position 24 staple with tt linker and anti-Quimby handle"
$arr[0,8] = "cc6hb h5 position 25 staple with tt linker and anti-Quimby handle"
$arr[0,9] = "cc6hb h5 position 26 staple with tt linker and anti-Quimby handle"
$arr[0,10] = "cc6hb h5 position 27 staple with tt linker and anti-Quimby handle"
$arr[0,11] = "cc6hb h5 position 28 staple with tt linker and anti-Quimby handle"
$arr[0,12] = "cc6hb h5 position 29 staple with tt linker and anti-Quimby handle"
$arr[0,13] = "cc6hb h5 position 30 staple with tt linker and anti-Quimby handle"
$arr[0,14] = "cc6hb h5 position 31 staple with tt linker and anti-Quimby handle"
$arr[0,15] = "cc6hb h5 position 32 staple with tt linker and anti-Quimby handle"
$ws.Range("B9:Q9").Value = $arr


$ws = $wb.Worksheets.Item("Names")
$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "antiQuimby_h5_pos1"
$arr[0,1] = "antiQuimby_h5_pos2"
$arr[0,2] = "antiQuimby_h5_pos3"
$arr[0,3] = "antiQuimby_h5_pos4"
$arr[0,4] = "antiQuimby_h5_pos5"
$arr[0,5] = "antiQuimby_h5_pos6"
$arr[0,6] = "antiQuimby_h5_pos7"
$arr[0,7] = "antiQuimby_h5_pos8"
$arr[0,8] = "antiQuimby_h5_pos9"
$arr[0,9] = "antiQuimby_h5_pos10"
$arr[0,10] = "antiQuimby_h5_pos11"
$arr[0,11] = "antiQuimby_h5_pos12"
$arr[0,12] = "antiQuimby_h5_pos13"
$arr[0,13] = "antiQuimby_h5_pos14"
$arr[0,14] = "antiQuimby_h5_pos15"
$arr[0,15] = "antiQuimby_h5_pos16"
$ws.Range("B8:Q8").Value = $arr

$arr = New-Object 'object[,]' 1,16
$arr[0,0] = "antiQuimby_h5_pos17"
$arr[0,1] = "antiQuimby_h5_pos18"
$arr[0,2] = "antiQuimby_h5_pos19"
$arr[0,3] = "antiQuimby_h5_pos20"
$arr[0,4] = "antiQuimby_h5_pos21"
$arr[0,5] = "antiQuimby_h5_pos22"
$arr[0,6] = "antiQuimby_h5_pos23"
$arr[0,7] = "antiQuimby_h5_pos24"
$arr[0,8] = "antiQuimby_h5_pos25"
$arr[0,9] = "antiQuimby_h5_pos26"
$arr[0,10] = "antiQuimby_h5_pos27"
$arr[0,11] = "antiQuimby_h5_pos28"
$arr[0,12] = "antiQuimby_h5_pos29"
$arr[0,13] = "antiQuimby_h5_pos30"
$arr[0,14] = "antiQuimby_h5_pos31"
$arr[0,15] = "antiQuimby_h5_pos32"
$ws.Range("B9:Q9").Value = $arr

